$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "1.000"); force text
# storage via NumberFormat so Excel does not coerce them to numbers, then
# reset the style back to Normal so no stray number-format style is left on
# the cell (matches the original file, where these cells carry no style).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.864.72'
$ws.Range("E2").Value = '  -2.71%  '

$ws.Range("D3").Value = '1.792.79'

$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '316.78'
$ws.Range("E5").Value = '  -0.16%  '

$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").Value = '0.5323'
$ws.Range("E7").Value = '  -0.37%  '

$ws.Range("D8").Value = '0.3863'
$ws.Range("E8").Value = '  +2.59%  '

$ws.Range("D9").Value = '0.07446'
$ws.Range("E9").Value = '  -0.87%  '

$ws.Range("D10").Value = '41.45'
$ws.Range("E10").Value = '  -2.25%  '

$ws.Range("D11").Value = '1.085'
$ws.Range("E11").Value = '  -2.72%  '

$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  -0.03%  '

$ws.Range("D13").Value = '6.179'
$ws.Range("E13").Value = '  +0.36%  '

$ws.Range("D14").Value = '7.456'
$ws.Range("E14").Value = '  +1.07%  '

$ws.Range("D15").Value = '20.30'
$ws.Range("E15").Value = '  -2.07%  '

$ws.Range("D16").Value = '1.797.39'
$ws.Range("E16").Value = '  +0.03%  '

$ws.Range("D17").Value = '88.27'
$ws.Range("E17").Value = '  -2.18%  '

$ws.Range("D18").Value = '0.00001059'
$ws.Range("E18").Value = '  -0.58%  '

$ws.Range("D19").Value = '0.06535'
$ws.Range("E19").Value = '  +1.17%  '

$ws.Range("E20").Value = '  -0.07%  '

$ws.Range("D21").Value = '17.22'
$ws.Range("E21").Value = '  -0.33%  '

$ws.Range("D22").Value = '5.969'
$ws.Range("E22").Value = '  +0.90%  '

$ws.Range("D23").Value = '27.900.00'
$ws.Range("E23").Value = '  -2.60%  '

$ws.Range("D24").Value = '11.11'
$ws.Range("E24").Value = '  -0.10%  '

$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("D26").Value = '156.28'
$ws.Range("E26").Value = '  -1.43%  '

$ws.Range("D27").Value = '20.13'
$ws.Range("E27").Value = '  -1.65%  '

$ws.Range("D28").Value = '2.000.97'
$ws.Range("E28").Value = '  -0.19%  '

$ws.Range("D29").Value = '2.305'
$ws.Range("E29").Value = '  -2.14%  '

$ws.Range("D30").Value = '121.58'
$ws.Range("E30").Value = '  -0.97%  '

$ws.Range("D31").Value = '0.1093'
$ws.Range("E31").Value = '  +2.61%  '

$ws.Range("D32").Value = '1.098'
$ws.Range("E32").Value = '  -0.79%  '

$ws.Range("D33").Value = '3.666'
$ws.Range("E33").Value = '  -0.36%  '

$ws.Range("D34").Value = '5.503'
$ws.Range("E34").Value = '  -2.75%  '

$ws.Range("D35").Value = '0.06930'
$ws.Range("E35").Value = '  +7.78%  '

$ws.Range("D36").Value = '0.2199'
$ws.Range("E36").Value = '  -2.37%  '

$ws.Range("D37").Value = '0.02268'
$ws.Range("E37").Value = '  -1.59%  '

$ws.Range("D38").Value = '5.054'
$ws.Range("E38").Value = '  +0.25%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '8.398'
$ws.Range("E39").Value = '  -4.20%  '

$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '11.24'
$ws.Range("E40").Value = '  -0.29%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '1.186'
$ws.Range("E41").Value = '  -1.04%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.6110'
$ws.Range("E42").Value = '  -1.90%  '

$ws.Range("D43").Value = '1.414'
$ws.Range("E43").Value = '  -0.63%  '

$ws.Range("D44").Value = '13.33'
$ws.Range("E44").Value = '  +0.76%  '

$ws.Range("D45").Value = '3.681'
$ws.Range("E45").Value = '  -0.26%  '

$ws.Range("D46").Value = '0.5707'
$ws.Range("E46").Value = '  -2.61%  '

$ws.Range("D47").Value = '125.05'
$ws.Range("E47").Value = '  -1.16%  '

$ws.Range("D48").Value = '1.914'
$ws.Range("E48").Value = '  -1.62%  '

$ws.Range("D49").Value = '1.175'
$ws.Range("E49").Value = '  +1.70%  '

$ws.Range("D50").Value = '0.06796'

$ws.Range("D51").Value = '71.31'
$ws.Range("E51").Value = '  -1.14%  '

$ws.Range("D2:D51").Style = "Normal"
